# Applies the numeric value updates (F = wanted-to-go count, G = min price)
# described in the commit diff, per worksheet.
# Sheet order in this workbook: 1=Exhibition, 2=Performance, 3=Local life, 4=All types
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Exhibition ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 417
$ws.Range("F3").Value = 1034
$ws.Range("F4").Value = 5814
$ws.Range("G4").Value = 75
$ws.Range("F6").Value = 1021
$ws.Range("F7").Value = 1030
$ws.Range("F8").Value = 843
$ws.Range("F9").Value = 85
$ws.Range("F11").Value = 612
$ws.Range("F12").Value = 43
$ws.Range("F15").Value = 1978
$ws.Range("F16").Value = 1507
$ws.Range("F17").Value = 1059
$ws.Range("F20").Value = 395
$ws.Range("F21").Value = 628
$ws.Range("F22").Value = 222
$ws.Range("F23").Value = 1068
$ws.Range("F26").Value = 3448
$ws.Range("F28").Value = 130
$ws.Range("F29").Value = 101
$ws.Range("F30").Value = 154
$ws.Range("F32").Value = 483
$ws.Range("F34").Value = 48
$ws.Range("F35").Value = 19
$ws.Range("F38").Value = 812
$ws.Range("F39").Value = 103
$ws.Range("F40").Value = 69
$ws.Range("F41").Value = 76
$ws.Range("F42").Value = 85

# --- Sheet 2: Performance ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 558
$ws.Range("F6").Value = 305

# --- Sheet 3: Local life ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 239

# --- Sheet 4: All types ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 239
$ws.Range("F3").Value = 417
$ws.Range("F4").Value = 1034
$ws.Range("F6").Value = 5814
$ws.Range("G6").Value = 75
$ws.Range("F8").Value = 1021
$ws.Range("F10").Value = 558
$ws.Range("F11").Value = 1030
$ws.Range("F12").Value = 843
$ws.Range("F14").Value = 305
$ws.Range("F15").Value = 85
$ws.Range("F17").Value = 612
$ws.Range("F18").Value = 44
$ws.Range("F22").Value = 1978
$ws.Range("F23").Value = 1507
$ws.Range("F24").Value = 1059
$ws.Range("F27").Value = 395
$ws.Range("F29").Value = 628
$ws.Range("F30").Value = 222
$ws.Range("F31").Value = 1068
$ws.Range("F32").Value = 3448
$ws.Range("F34").Value = 130
$ws.Range("F35").Value = 101
$ws.Range("F36").Value = 154
$ws.Range("F38").Value = 483
$ws.Range("F40").Value = 48
$ws.Range("F41").Value = 19
$ws.Range("F43").Value = 812
$ws.Range("F44").Value = 103
$ws.Range("F45").Value = 76
$ws.Range("F46").Value = 85

